$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.782.81'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '2.190.85'
$ws.Range('E3').Value = '  -2.16%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '292.76'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '86.55'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').Value = '  -1.62%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.466'
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0774'
$ws.Range('E10').Value = '  -2.27%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '29.76'
$ws.Range('E11').Value = '  -3.55%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '49.57'
$ws.Range('E12').Value = '  +5.49%  '
$ws.Range('E13').Value = '  +2.21%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.38'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').Value = '2.536.81'
$ws.Range('E15').Value = '  -1.40%  '
$ws.Range('D16').Value = '2.257.75'
$ws.Range('E16').Value = '  +2.71%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.59'
$ws.Range('E17').Value = '  -3.80%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.719'
$ws.Range('E18').Value = '  -1.64%  '
$ws.Range('D19').Value = '39.716.38'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('D20').Value = '0.0₃0878'
$ws.Range('E20').Value = '  -1.48%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.13'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.71'
$ws.Range('E22').Value = '  -1.78%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '64.87'
$ws.Range('E23').Value = '  -0.71%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '235.74'
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.44'
$ws.Range('E26').Value = '  -1.21%  '
$ws.Range('E27').Value = '  -2.96%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '22.31'
$ws.Range('E28').Value = '  -2.34%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.14'
$ws.Range('E29').Value = '  -3.64%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.09'
$ws.Range('E30').Value = '  -2.45%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '156.39'
$ws.Range('E31').Value = '  +2.52%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '31.04'
$ws.Range('E32').Value = '  -6.69%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.84'
$ws.Range('E34').Value = '  -2.20%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0701'
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('E36').Value = '  -2.37%  '
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.112'
$ws.Range('E37').Value = '  +0.46%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.78'
$ws.Range('E38').Value = '  -0.64%  '
$ws.Range('E39').Value = '  -3.46%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '15.10'
$ws.Range('E40').Value = '  -7.00%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.66'
$ws.Range('E41').Value = '  -3.42%  '
$ws.Range('D42').Value = '2.104.25'
$ws.Range('E42').Value = '  +3.09%  '
$ws.Range('E43').Value = '  -3.51%  '
$ws.Range('E44').Value = '  -2.05%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0265'
$ws.Range('E45').Value = '  -2.18%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '17.19'
$ws.Range('E46').Value = '  +1.10%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.56'
$ws.Range('E47').Value = '  -5.33%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.62'
$ws.Range('E48').Value = '  +1.31%  '
$ws.Range('D49').Value = '2.411.17'
$ws.Range('E49').Value = '  -2.03%  '
$ws.Range('E50').Value = '  +3.09%  '
$ws.Range('E51').Value = '  +0.78%  '
